$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above row 11 ("Description") to make room for "Jurisdiction"
$ws.Rows.Item(11).Insert()

# Populate the new "Jurisdiction" row
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update the URL value (row 2, column B)
$ws.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/CodeSystem/IntervalReason"

# Update the Date value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"
